$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bolivia (Plurinational State of) -> Bolivia, rows 126-156
for ($r = 126; $r -le 156; $r++) {
    $ws.Cells.Item($r, 1).Value = "Bolivia"
}

# United States of America (the) -> United States, rows 994-1024
for ($r = 994; $r -le 1024; $r++) {
    $ws.Cells.Item($r, 1).Value = "United States"
}

# Venezuela (Bolivarian Republic of) -> Venezuela, rows 1056-1086
for ($r = 1056; $r -le 1086; $r++) {
    $ws.Cells.Item($r, 1).Value = "Venezuela"
}
